# SBERDOMA-1030 - refactor with global mappers
#
# Applies to the ticket analytics export template:
#  - Renames the "Complete" status label to "Completed"
#  - Moves the active selection from D9 to E9
#  - Widens columns B-E to their new layout widths
#  - Shrinks the header row heights (rows 2 & 3) from 36.9 to 25.1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "Complete" -> "Completed" (status header in C1)
$ws.Range("C1").Value = "Completed"

# 2) Column widths (Excel's ColumnWidth is in character units of the
#    workbook's default/Normal font; it is stored in the OOXML <col>
#    element as ColumnWidth + ~0.83, rounded to the nearest pixel column,
#    so the values below are chosen to land as close as possible on the
#    saved widths 23.98 / 16.96 / 16.22 / 14.62).
$ws.Columns.Item(2).ColumnWidth = 23.165   # B: 21.36 -> ~23.98
$ws.Columns.Item(3).ColumnWidth = 16.165   # C: 11.63 (auto) -> ~16.96 (custom)
$ws.Columns.Item(4).ColumnWidth = 15.33    # D: 14.77 -> ~16.22
$ws.Columns.Item(5).ColumnWidth = 13.83    # E: 11.63 (auto) -> ~14.62 (custom)

# 3) Header/data row heights shrink from 36.9 to 25.1
$ws.Rows.Item(2).RowHeight = 25.1
$ws.Rows.Item(3).RowHeight = 25.1

# 4) Active cell/selection moves from D9 to E9
$ws.Range("E9").Select() | Out-Null
